$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row for "Edinaldo De Paiva Santos" (row 4), which
# shifts subsequent rows (Filipe, Jordan, Lucas) up by one.
$ws.Rows.Item(4).Delete()

# Restore the selection to match the post-delete state (Excel leaves the
# whole row selected after a "Delete Row" operation).
$ws.Range("A4:XFD4").Select()
